# Updated vignette tables (labels.xlsx):
# Replace the old Storability/good/bad table with the new
# Samples / Class / Colorcode table (Sample1..Sample8, control / treatment1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "Samples"
$ws.Range("B1").Value = "Class"
$ws.Range("C1").Value = "Colorcode"

# control group (Sample1-4)
$ws.Range("A2").Value = "Sample1"
$ws.Range("B2").Value = "control"
$ws.Range("C2").Value = "#b3ca81"

$ws.Range("A3").Value = "Sample2"
$ws.Range("B3").Value = "control"
$ws.Range("C3").Value = "#b3ca82"

$ws.Range("A4").Value = "Sample3"
$ws.Range("B4").Value = "control"
$ws.Range("C4").Value = "#b3ca83"

$ws.Range("A5").Value = "Sample4"
$ws.Range("B5").Value = "control"
$ws.Range("C5").Value = "#b3ca84"

# treatment1 group (Sample5-8)
$ws.Range("A6").Value = "Sample5"
$ws.Range("B6").Value = "treatment1"
$ws.Range("C6").Value = "red"

$ws.Range("A7").Value = "Sample6"
$ws.Range("B7").Value = "treatment1"
$ws.Range("C7").Value = "red"

$ws.Range("A8").Value = "Sample7"
$ws.Range("B8").Value = "treatment1"
$ws.Range("C8").Value = "red"

$ws.Range("A9").Value = "Sample8"
$ws.Range("B9").Value = "treatment1"
$ws.Range("C9").Value = "red"

# Move the active selection to H8, matching the saved workbook state.
$ws.Range("H8").Select()
